$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Vaccine import columns -------------------------------------------------
# The "Administration Date" columns (Vaccine 1 -> DA, Vaccine 2 -> DF) must be
# stored as literal text (e.g. "2020-06-01") rather than parsed into Excel
# date serials, so format them as text before writing any values into them.
$ws.Range("DA1:DA7").NumberFormat = "@"
$ws.Range("DF1:DF4").NumberFormat = "@"

# Header row (row 1): Vaccine 1 (CY:DC) and Vaccine 2 (DD:DH) field labels.
$ws.Range("CY1").Value = "Vaccine 1 Group Name"
$ws.Range("CZ1").Value = "Vaccine 1 Product Name"
$ws.Range("DA1").Value = "Vaccine 1 Administration Date"
$ws.Range("DB1").Value = "Vaccine 1 Dose Number"
$ws.Range("DC1").Value = "Vaccine 1 Notes"
$ws.Range("DD1").Value = "Vaccine 2 Group Name"
$ws.Range("DE1").Value = "Vaccine 2 Product Name"
$ws.Range("DF1").Value = "Vaccine 2 Administration Date"
$ws.Range("DG1").Value = "Vaccine 2 Dose Number"
$ws.Range("DH1").Value = "Vaccine 2 Notes"

# Row 2
$ws.Range("CY2").Value = "COVID-19"
$ws.Range("CZ2").Value = "Moderna COVID-19 Vaccine"
$ws.Range("DA2").Value = "2020-06-01"
$ws.Range("DB2").Value = 1
$ws.Range("DC2").Value = "notes 1"
$ws.Range("DD2").Value = "COVID-19"
$ws.Range("DE2").Value = "Moderna COVID-19 Vaccine"
$ws.Range("DF2").Value = "2020-06-20"
$ws.Range("DG2").Value = 2
$ws.Range("DH2").Value = "notes 2"

# Row 3
$ws.Range("CY3").Value = "COVID-19"
$ws.Range("CZ3").Value = "Pfizer-BioNTech COVID-19 Vaccine"
$ws.Range("DA3").Value = "2020-06-02"
$ws.Range("DB3").Value = 1
$ws.Range("DD3").Value = "COVID-19"
$ws.Range("DE3").Value = "Pfizer-BioNTech COVID-19 Vaccine"
$ws.Range("DF3").Value = "2020-06-21"
$ws.Range("DG3").Value = 2

# Row 4
$ws.Range("CY4").Value = "COVID-19"
$ws.Range("CZ4").Value = "Unknown"
$ws.Range("DA4").Value = "2020-06-04"
$ws.Range("DB4").Value = 1
$ws.Range("DD4").Value = "COVID-19"
$ws.Range("DE4").Value = "Unknown"
$ws.Range("DF4").Value = "2020-06-22"
$ws.Range("DG4").Value = 2

# Row 5
$ws.Range("CY5").Value = "COVID-19"
$ws.Range("CZ5").Value = "Moderna COVID-19 Vaccine"
$ws.Range("DA5").Value = "2020-06-01"
$ws.Range("DB5").Value = 1

# Row 6
$ws.Range("CY6").Value = "COVID-19"
$ws.Range("CZ6").Value = "Janssen (J&J) COVID-19 Vaccine"
$ws.Range("DA6").Value = "2020-06-03"
$ws.Range("DB6").Value = 1

# Row 7
$ws.Range("CY7").Value = "COVID-19"
$ws.Range("CZ7").Value = "Unknown"
$ws.Range("DA7").Value = "2020-06-02"
$ws.Range("DB7").Value = 1

# --- Column widths -----------------------------------------------------------
# Approximate the auto-fit widths Excel computed for the new columns (103-112 /
# CY-DH) when it saved the workbook after this edit. This engine stores
# ColumnWidth with a constant +5/7 offset versus the persisted <col width=.../>
# value, so back that offset out of the target widths we want to land on.
$targetWidths = @{
  103 = 20.33203125
  104 = 31
  105 = 25.6640625
  106 = 21.1640625
  107 = 14.5
  108 = 20.33203125
  109 = 31
  110 = 25.6640625
  111 = 21.1640625
  112 = 14.5
}
foreach ($col in $targetWidths.Keys) {
  $ws.Columns.Item($col).ColumnWidth = $targetWidths[$col] - (5/7)
}

# Reset the active selection back to A1 (matches the saved view after the
# edit, instead of leaving the stale CU9 selection from before the import).
[void]$ws.Range("A1").Select()
